$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Flight Mission Cycle": the summary table shrinks from 5 settings rows
#    down to just "Writing" (now with Duration 0); the Piano / Light switch /
#    Typing rows move out into their own dedicated sheets (see below).
# ---------------------------------------------------------------------------
$fmc = $wb.Worksheets.Item("Flight Mission Cycle")
$fmc.Rows("3:5").Delete()
$fmc.Range("B2").Value = 0

# ---------------------------------------------------------------------------
# 2. Insert three new sheets between "Flight Mission Cycle" and "Writing".
#    Each new sheet is inserted immediately after "Flight Mission Cycle", so
#    inserting Piano, then Light switch, then Typing (in that order) leaves
#    them in final left-to-right order: Typing, Light switch, Piano.
#    (The worksheet reference returned by Add() can go stale once another
#    sheet is added, so each anchor/handle is re-fetched fresh by name.)
# ---------------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("Flight Mission Cycle")
$piano = $wb.Worksheets.Add($null, $anchor)
$piano.Name = "Piano"

$anchor = $wb.Worksheets.Item("Flight Mission Cycle")
$lightswitch = $wb.Worksheets.Add($null, $anchor)
$lightswitch.Name = "Light switch"

$anchor = $wb.Worksheets.Item("Flight Mission Cycle")
$typing = $wb.Worksheets.Add($null, $anchor)
$typing.Name = "Typing"

# ---------------------------------------------------------------------------
# 3. Populate "Typing" data (re-fetched by name to avoid any stale handle).
# ---------------------------------------------------------------------------
$typing = $wb.Worksheets.Item("Typing")
$typing.Range("A1").Value = "Typing"
$typing.Range("A2").Value = "Force"
$typing.Range("B2").Value = 0
$typing.Range("C2").Value = 10
$typing.Range("D2").Value = 0
$typing.Range("A3").Value = "Duration"
$typing.Range("B3").Value = 10
$typing.Range("C3").Value = 30
$typing.Range("D3").Value = 40
$typing.Range("A4").Value = "Max_RoM"
$typing.Range("B4").Value = -30
$typing.Range("A5").Value = "Min_RoM"
$typing.Range("B5").Value = 60
$typing.Range("A6").Value = "Period"
$typing.Range("B6").Value = 0

# ---------------------------------------------------------------------------
# 4. Populate "Light switch" data.
# ---------------------------------------------------------------------------
$lightswitch = $wb.Worksheets.Item("Light switch")
$lightswitch.Range("A1").Value = "Piano"
$lightswitch.Range("A2").Value = "Force"
$lightswitch.Range("A3").Value = "Duration"
$lightswitch.Range("A4").Value = "Max_RoM"
$lightswitch.Range("B4").Value = -30
$lightswitch.Range("A5").Value = "Min_RoM"
$lightswitch.Range("B5").Value = 60
$lightswitch.Range("A6").Value = "Period"
$lightswitch.Range("B6").Value = 10

# ---------------------------------------------------------------------------
# 5. Populate "Piano" data.
# ---------------------------------------------------------------------------
$piano = $wb.Worksheets.Item("Piano")
$piano.Range("A1").Value = "Piano"
$piano.Range("A2").Value = "Force"
$piano.Range("B2").Value = 0
$piano.Range("C2").Value = 50
$piano.Range("D2").Value = 50
$piano.Range("E2").Value = 30
$piano.Range("F2").Value = 30
$piano.Range("G2").Value = 0
$piano.Range("A3").Value = "Duration"
$piano.Range("B3").Value = 0
$piano.Range("C3").Value = 20
$piano.Range("D3").Value = 30
$piano.Range("E3").Value = 10
$piano.Range("F3").Value = 100
$piano.Range("G3").Value = 10
$piano.Range("A4").Value = "Max_RoM"
$piano.Range("B4").Value = 0
$piano.Range("A5").Value = "Min_RoM"
$piano.Range("B5").Value = -60
$piano.Range("A6").Value = "Period"
$piano.Range("B6").Value = 10

# ---------------------------------------------------------------------------
# 6. View-state bits: "Flight Mission Cycle" becomes the active tab/sheet
#    (activeTab index 1, 0-based), with the selected cell moved to H17.
#    "Writing" loses tabSelected and gets a whole-column selection instead.
#    (Set the non-active sheet's selection first, then activate/select the
#    sheet that should end up as the active tab last.)
# ---------------------------------------------------------------------------
$writing = $wb.Worksheets.Item("Writing")
$writing.Range("A1:A6").Select()

$fmc = $wb.Worksheets.Item("Flight Mission Cycle")
$fmc.Activate()
$fmc.Range("H17").Select()
